# "aggiornamento fino a 6/03" - append 3 new daily rows (r=245..247) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44319, 0, 34, 193.5226819966987),
    @(44320, 1, 23, 130.9124025271786),
    @(44321, 0, 22, 125.2205589390404)
)

$row = 245
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# The date column (A) uses a dedicated date style (same as every other row in
# the column, e.g. A244) - copy that formatting onto the new date cells.
$ws.Range("A244").Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)
